$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 9845.454
$ws.Range("I33").Value = 100
$ws.Range("K33").Value = 100
$ws.Range("M33").Value = 129

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3149.3333
$ws.Range("I113").Value = 3161.3333
$ws.Range("J113").Value = 3143.3333
$ws.Range("K113").Value = 3161.3333
$ws.Range("L113").Value = 3143.3333
$ws.Range("M113").Value = 92.66670000000022
$ws.Range("N113").Value = -9651.3333

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 936.2162
$ws.Range("I129").Value = 299
$ws.Range("K129").Value = 897
$ws.Range("M129").Value = 4103

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2509.6216
$ws.Range("I138").Value = 1254.9434
$ws.Range("J138").Value = 5676.1904
$ws.Range("K138").Value = 3764.8302
$ws.Range("L138").Value = 17028.5712
$ws.Range("M138").Value = 1375.1698
$ws.Range("N138").Value = -27308.5712

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 399965.12
$ws.Range("I32").Value = 2640.4285
$ws.Range("J32").Value = 3181238
$ws.Range("K32").Value = 2640.4285
$ws.Range("L32").Value = 3181238
$ws.Range("M32").Value = -2353.4285
$ws.Range("N32").Value = -3181812

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1117.8611
$ws.Range("I61").Value = 733.6786
$ws.Range("J61").Value = 2462.5
$ws.Range("K61").Value = 733.6786
$ws.Range("L61").Value = 2462.5
$ws.Range("M61").Value = -521.6786
$ws.Range("N61").Value = -2886.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1496
$ws.Range("I132").Value = 1018.75
$ws.Range("J132").Value = 5314
$ws.Range("K132").Value = 3056.25
$ws.Range("L132").Value = 15942
$ws.Range("M132").Value = -526.25
$ws.Range("N132").Value = -21002

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1117.8611
$ws.Range("I136").Value = 733.6786
$ws.Range("J136").Value = 2462.5
$ws.Range("K136").Value = 2201.0358
$ws.Range("L136").Value = 7387.5
$ws.Range("M136").Value = 348.9642000000003
$ws.Range("N136").Value = -12487.5

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 46075
$ws.Range("J139").Value = 46075
$ws.Range("L139").Value = 46075
$ws.Range("N139").Value = -56355

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 36254.332
$ws.Range("I81").Value = 20709
$ws.Range("J81").Value = 38197.5
$ws.Range("K81").Value = 20709
$ws.Range("L81").Value = 38197.5
$ws.Range("M81").Value = -19648
$ws.Range("N81").Value = -40319.5

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 36254.332
$ws.Range("I84").Value = 20709
$ws.Range("J84").Value = 38197.5
$ws.Range("K84").Value = 62127
$ws.Range("L84").Value = 114592.5
$ws.Range("M84").Value = -56823
$ws.Range("N84").Value = -125200.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 26237.6
$ws.Range("I132").Value = 1170.5186
$ws.Range("J132").Value = 78300
$ws.Range("K132").Value = 3511.5558
$ws.Range("L132").Value = 234900
$ws.Range("M132").Value = -981.5558000000001
$ws.Range("N132").Value = -239960

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1158.0571
$ws.Range("I134").Value = 1201.1428
$ws.Range("J134").Value = 985.7143
$ws.Range("K134").Value = 3603.4284
$ws.Range("L134").Value = 2957.1429
$ws.Range("M134").Value = -1068.4284
$ws.Range("N134").Value = -8027.1429

# CUL row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1561.2667
$ws.Range("I116").Value = 809.6667
$ws.Range("J116").Value = 1749.1666
$ws.Range("K116").Value = 2429.0001
$ws.Range("L116").Value = 5247.4998
$ws.Range("M116").Value = 1012.9999
$ws.Range("N116").Value = -12131.4998

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17858066
$ws.Range("I131").Value = 1056.6666
$ws.Range("J131").Value = 26316650
$ws.Range("K131").Value = 3169.9998
$ws.Range("L131").Value = 78949950
$ws.Range("M131").Value = 1870.0002
$ws.Range("N131").Value = -78960030

# GSM row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 51997.5
$ws.Range("J133").Value = 51997.5
$ws.Range("L133").Value = 51997.5
$ws.Range("N133").Value = -62117.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2533.1667
$ws.Range("I122").Value = 2490.375
$ws.Range("J122").Value = 2618.75
$ws.Range("K122").Value = 7471.125
$ws.Range("L122").Value = 7856.25
$ws.Range("M122").Value = -5021.125
$ws.Range("N122").Value = -12756.25

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# WVR row 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -893

# WVR row 33
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

# WVR row 36
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

# WVR row 40
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 6000
$ws.Range("J40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("N40").Value = -6298

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18118364
$ws.Range("I132").Value = 27174932
$ws.Range("J132").Value = 5227.826
$ws.Range("K132").Value = 81524796
$ws.Range("L132").Value = 15683.478
$ws.Range("M132").Value = -81522266
$ws.Range("N132").Value = -20743.478

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3121.1
$ws.Range("I136").Value = 4154.407
$ws.Range("J136").Value = 975
$ws.Range("K136").Value = 12463.221
$ws.Range("L136").Value = 2925
$ws.Range("M136").Value = -9913.221000000001
$ws.Range("N136").Value = -8025
